$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet/tab from SCD0335 to SCD0025
$ws.Name = "SCD0025"

# Update the data row (row 2): TC_ID, TEST_SCENARIO_DESC, EXPECTED_RESULT
# B2 TC_ID: DGS-350 -> SCD0025-005
$ws.Range("B2").Value = "SCD0025-005"
# C2 TEST_SCENARIO_DESC
$ws.Range("C2").Value = "Normal Skenario Admin SLN mengakses Report Log Pengiriman WA pada Digisales Portal"
# E2 EXPECTED_RESULT
$ws.Range("E2").Value = "Berhasil memunculkan data report"

# Column B needs to widen (bestFit) to fit the longer TC_ID text
$ws.Columns.Item(2).AutoFit() | Out-Null

# Move the active selection to B3 (also resets the scrolled view back to column A)
$ws.Range("B3").Select()
